$d = $word.ActiveDocument

$rng = $d.Content
$rng.Find.Execute("Visual Studio (C#)", $true, $false, $false, $false, $false, $true, 1, $false, "Visual Studio Community 2013 (C#)", 2)
